# Updated symbol list with GitHub Actions — refresh scraped coin prices /
# symbols on Sheet1 (cryptos.xlsx export).
#
# All target cells in this sheet are stored as TEXT (even the numeric-looking
# price column D), so each cell is explicitly formatted as Text ("@") before
# its value is assigned — otherwise a numeric-looking string like "249.23"
# assigned to a General-formatted cell gets auto-coerced into a real number,
# same as typing it straight into Excel would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @("D2", "249.23"),
    @("D3", "21.66"),
    @("D4", "5.547"),
    @("D5", "0.05660"),
    @("D6", "6.453"),
    @("D7", "0.8000"),
    @("D8", "1.045"),
    @("D9", "0.1439"),
    @("D10", "0.07324"),
    @("D11", "0.03132"),
    @("D12", "0.02916"),
    @("D14", "0.001664"),
    @("D15", "3.217"),
    @("D16", "0.04731"),
    @("D17", "0.0005820"),
    @("E17", "16OneONE"),
    @("D18", "0.006451"),
    @("D19", "0.005077"),
    @("D20", "0.001050"),
    @("D22", "3.978"),
    @("D23", "3.374"),
    @("D25", "0.3269"),
    @("D27", "0.0003200"),
    @("D40", "0.04143"),
    @("B41", "BKEXToken"),
    @("C41", "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"),
    @("D41", "0.1045"),
    @("E41", "40BKEXTokenBKK"),
    @("D42", "0.002972"),
    @("B43", "KickToken"),
    @("C43", "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"),
    @("D43", "0.006944"),
    @("E43", "42KickTokenKICK"),
    @("D44", "0.009405"),
    @("D45", "0.00005643"),
    @("D47", "0.7852"),
    @("D48", "0.01658"),
    @("E48", "47BOLOBOLOWorstin24h"),
    @("D49", "0.00002100"),
    @("D50", "0.01010")
)

foreach ($item in $changes) {
    $addr = $item[0]
    $text = $item[1]
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    # Drop the temporary Text number-format again so the cell's style index
    # matches the untouched cells around it (only its stored value/type
    # changed) — ClearFormats() resets NumberFormat back to General/style 0
    # without touching the value we just wrote.
    $cell.ClearFormats()
}
